{"js": "const replacements = [\n  [\"85-30=55\", \"98-88=10\"],\n  [\"89+3=92\", \"13+36=49\"],\n  [\"50+27=77\", \"21+6=27\"],\n  [\"67+22=89\", \"37-31=6\"],\n  [\"58-25=33\", \"42+1=43\"],\n  [\"80-24=56\", \"4+72=76\"],\n  [\"44-25=19\", \"42-17=25\"],\n  [\"67-2=65\", \"96-16=80\"],\n  [\"63-44=19\", \"45+0=45\"],\n  [\"62+18=80\", \"67+29=96\"],\n  [\"52-20=32\", \"88-17=71\"],\n  [\"0+28=28\", \"77+11=88\"],\n  [\"28+63=91\", \"66-9=57\"],\n  [\"97-44=53\", \"83-49=34\"],\n  [\"48-46=2\", \"3+4=7\"],\n  [\"64+9=73\", \"44+12=56\"],\n  [\"48+6=54\", \"5+36=41\"],\n  [\"38+3=41\", \"63-62=1\"],\n  [\"3+26=29\", \"94-76=18\"],\n  [\"97-83=14\", \"48+40=88\"],\n  [\"5+23=28\", \"45+27=72\"],\n  [\"86-47=39\", \"40+48=88\"],\n  [\"10+74=84\", \"93-64=29\"],\n  [\"66-41=25\", \"79-16=63\"],\n  [\"1+88=89\", \"72-24=48\"],\n  [\"69-35=34\", \"30-15=15\"],\n  [\"31-20=11\", \"22+76=98\"],\n  [\"10+53=63\", \"21+22=43\"],\n  [\"36-25=11\", \"53+29=82\"],\n  [\"47+41=88\", \"21+49=70\"],\n  [\"24-4=20\", \"80-65=15\"],\n  [\"12+29=41\", \"62-6=56\"],\n  [\"60-8=52\", \"61+12=73\"],\n  [\"83-6=77\", \"6+27=33\"],\n  [\"41-3=38\", \"33+0=33\"],\n  [\"75+14=89\", \"19+1=20\"],\n  [\"6+85=91\", \"17+28=45\"],\n  [\"7+82=89\", \"4-0=4\"],\n  [\"57-5=52\", \"22+35=57\"],\n  [\"88-0=88\", \"98-0=98\"],\n  [\"44-38=6\", \"21+51=72\"],\n  [\"33+44=77\", \"92-74=18\"],\n  [\"17+72=89\", \"20+41=61\"],\n  [\"9+62=71\", \"10+4=14\"],\n  [\"96-96=0\", \"72+6=78\"],\n  [\"64+23=87\", \"50+23=73\"],\n  [\"97-51=46\", \"47+46=93\"],\n  [\"0+84=84\", \"54-46=8\"],\n  [\"0+39=39\", \"67+19=86\"],\n  [\"71-42=29\", \"79-33=46\"],\n  [\"81+0=81\", \"85-5=80\"],\n  [\"30+47=77\", \"73-31=42\"],\n  [\"15+80=95\", \"11+36=47\"],\n  [\"65-34=31\", \"8+84=92\"],\n  [\"48+34=82\", \"64-59=5\"],\n  [\"72-64=8\", \"12-3=9\"],\n  [\"58-18=40\", \"45-17=28\"],\n  [\"13+83=96\", \"34+12=46\"],\n  [\"63+34=97\", \"68-14=54\"],\n  [\"29+21=50\", \"57-7=50\"],\n  [\"62-12=50\", \"77-77=0\"],\n  [\"53+17=70\", \"29+36=65\"],\n  [\"23+10=33\", \"30+61=91\"],\n  [\"42+29=71\", \"52+25=77\"],\n  [\"82-48=34\", \"50+28=78\"],\n  [\"53-41=12\", \"68-19=49\"],\n  [\"25+49=74\", \"14+32=46\"],\n  [\"93-9=84\", \"94-23=71\"],\n  [\"88-43=45\", \"23-21=2\"],\n  [\"45-21=24\", \"74+12=86\"],\n  [\"79-51=28\", \"92-59=33\"],\n  [\"92-52=40\", \"18+21=39\"],\n  [\"2+46=48\", \"71-47=24\"],\n  [\"83+4=87\", \"84+14=98\"],\n  [\"24+60=84\", \"72-15=57\"],\n  [\"48+49=97\", \"74-1=73\"],\n  [\"36+3=39\", \"31+40=71\"],\n  [\"43-32=11\", \"84-27=57\"],\n  [\"87-2=85\", \"67+3=70\"],\n  [\"89-40=49\", \"18+11=29\"],\n  [\"18+32=50\", \"64-42=22\"],\n  [\"35+50=85\", \"92-1=91\"],\n  [\"85-60=25\", \"63-9=54\"],\n  [\"64+1=65\", \"91+4=95\"],\n  [\"18+80=98\", \"89-84=5\"],\n  [\"69+25=94\", \"82-12=70\"],\n  [\"26+13=39\", \"32-7=25\"],\n  [\"42-33=9\", \"45+4=49\"],\n  [\"25-18=7\", \"86-66=20\"],\n  [\"80-49=31\", \"83-69=14\"],\n  [\"93-6=87\", \"46+40=86\"],\n  [\"93-17=76\", \"38+53=91\"],\n  [\"46-17=29\", \"44+51=95\"],\n  [\"52+2=54\", \"90+3=93\"],\n  [\"98-47=51\", \"8+7=15\"],\n  [\"62+30=92\", \"88-82=6\"],\n  [\"31-27=4\", \"76-15=61\"],\n  [\"26+67=93\", \"25-22=3\"],\n  [\"38+22=60\", \"42+24=66\"],\n  [\"53+9=62\", \"74+9=83\"],\n];\n\nfor (const [oldText, newText] of replacements) {\n  const results = context.document.body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n  if (results.items.length === 0) {\n    throw new Error(\"No match found for: \" + oldText);\n  }\n  results.items[0].insertText(newText, Word.InsertLocation.replace);\n}\nawait context.sync();", "ps1": "$d = $word.ActiveDocument\n\n$replacements = @(\n  ,@(\"85-30=55\", \"98-88=10\")\n  ,@(\"89+3=92\", \"13+36=49\")\n  ,@(\"50+27=77\", \"21+6=27\")\n  ,@(\"67+22=89\", \"37-31=6\")\n  ,@(\"58-25=33\", \"42+1=43\")\n  ,@(\"80-24=56\", \"4+72=76\")\n  ,@(\"44-25=19\", \"42-17=25\")\n  ,@(\"67-2=65\", \"96-16=80\")\n  ,@(\"63-44=19\", \"45+0=45\")\n  ,@(\"62+18=80\", \"67+29=96\")\n  ,@(\"52-20=32\", \"88-17=71\")\n  ,@(\"0+28=28\", \"77+11=88\")\n  ,@(\"28+63=91\", \"66-9=57\")\n  ,@(\"97-44=53\", \"83-49=34\")\n  ,@(\"48-46=2\", \"3+4=7\")\n  ,@(\"64+9=73\", \"44+12=56\")\n  ,@(\"48+6=54\", \"5+36=41\")\n  ,@(\"38+3=41\", \"63-62=1\")\n  ,@(\"3+26=29\", \"94-76=18\")\n  ,@(\"97-83=14\", \"48+40=88\")\n  ,@(\"5+23=28\", \"45+27=72\")\n  ,@(\"86-47=39\", \"40+48=88\")\n  ,@(\"10+74=84\", \"93-64=29\")\n  ,@(\"66-41=25\", \"79-16=63\")\n  ,@(\"1+88=89\", \"72-24=48\")\n  ,@(\"69-35=34\", \"30-15=15\")\n  ,@(\"31-20=11\", \"22+76=98\")\n  ,@(\"10+53=63\", \"21+22=43\")\n  ,@(\"36-25=11\", \"53+29=82\")\n  ,@(\"47+41=88\", \"21+49=70\")\n  ,@(\"24-4=20\", \"80-65=15\")\n  ,@(\"12+29=41\", \"62-6=56\")\n  ,@(\"60-8=52\", \"61+12=73\")\n  ,@(\"83-6=77\", \"6+27=33\")\n  ,@(\"41-3=38\", \"33+0=33\")\n  ,@(\"75+14=89\", \"19+1=20\")\n  ,@(\"6+85=91\", \"17+28=45\")\n  ,@(\"7+82=89\", \"4-0=4\")\n  ,@(\"57-5=52\", \"22+35=57\")\n  ,@(\"88-0=88\", \"98-0=98\")\n  ,@(\"44-38=6\", \"21+51=72\")\n  ,@(\"33+44=77\", \"92-74=18\")\n  ,@(\"17+72=89\", \"20+41=61\")\n  ,@(\"9+62=71\", \"10+4=14\")\n  ,@(\"96-96=0\", \"72+6=78\")\n  ,@(\"64+23=87\", \"50+23=73\")\n  ,@(\"97-51=46\", \"47+46=93\")\n  ,@(\"0+84=84\", \"54-46=8\")\n  ,@(\"0+39=39\", \"67+19=86\")\n  ,@(\"71-42=29\", \"79-33=46\")\n  ,@(\"81+0=81\", \"85-5=80\")\n  ,@(\"30+47=77\", \"73-31=42\")\n  ,@(\"15+80=95\", \"11+36=47\")\n  ,@(\"65-34=31\", \"8+84=92\")\n  ,@(\"48+34=82\", \"64-59=5\")\n  ,@(\"72-64=8\", \"12-3=9\")\n  ,@(\"58-18=40\", \"45-17=28\")\n  ,@(\"13+83=96\", \"34+12=46\")\n  ,@(\"63+34=97\", \"68-14=54\")\n  ,@(\"29+21=50\", \"57-7=50\")\n  ,@(\"62-12=50\", \"77-77=0\")\n  ,@(\"53+17=70\", \"29+36=65\")\n  ,@(\"23+10=33\", \"30+61=91\")\n  ,@(\"42+29=71\", \"52+25=77\")\n  ,@(\"82-48=34\", \"50+28=78\")\n  ,@(\"53-41=12\", \"68-19=49\")\n  ,@(\"25+49=74\", \"14+32=46\")\n  ,@(\"93-9=84\", \"94-23=71\")\n  ,@(\"88-43=45\", \"23-21=2\")\n  ,@(\"45-21=24\", \"74+12=86\")\n  ,@(\"79-51=28\", \"92-59=33\")\n  ,@(\"92-52=40\", \"18+21=39\")\n  ,@(\"2+46=48\", \"71-47=24\")\n  ,@(\"83+4=87\", \"84+14=98\")\n  ,@(\"24+60=84\", \"72-15=57\")\n  ,@(\"48+49=97\", \"74-1=73\")\n  ,@(\"36+3=39\", \"31+40=71\")\n  ,@(\"43-32=11\", \"84-27=57\")\n  ,@(\"87-2=85\", \"67+3=70\")\n  ,@(\"89-40=49\", \"18+11=29\")\n  ,@(\"18+32=50\", \"64-42=22\")\n  ,@(\"35+50=85\", \"92-1=91\")\n  ,@(\"85-60=25\", \"63-9=54\")\n  ,@(\"64+1=65\", \"91+4=95\")\n  ,@(\"18+80=98\", \"89-84=5\")\n  ,@(\"69+25=94\", \"82-12=70\")\n  ,@(\"26+13=39\", \"32-7=25\")\n  ,@(\"42-33=9\", \"45+4=49\")\n  ,@(\"25-18=7\", \"86-66=20\")\n  ,@(\"80-49=31\", \"83-69=14\")\n  ,@(\"93-6=87\", \"46+40=86\")\n  ,@(\"93-17=76\", \"38+53=91\")\n  ,@(\"46-17=29\", \"44+51=95\")\n  ,@(\"52+2=54\", \"90+3=93\")\n  ,@(\"98-47=51\", \"8+7=15\")\n  ,@(\"62+30=92\", \"88-82=6\")\n  ,@(\"31-27=4\", \"76-15=61\")\n  ,@(\"26+67=93\", \"25-22=3\")\n  ,@(\"38+22=60\", \"42+24=66\")\n  ,@(\"53+9=62\", \"74+9=83\")\n)\n\nforeach ($pair in $replacements) {\n  $oldText = $pair[0]\n  $newText = $pair[1]\n  $find = $d.Content.Find\n  $find.ClearFormatting()\n  $find.Replacement.ClearFormatting()\n  $found = $find.Execute($oldText, $false, $false, $false, $false, $false, $true, 1, $false, $newText, 2)\n  if (-not $found) {\n    throw \"No match found for: $oldText\"\n  }\n}\n\nWrite-Output \"done\""}
